$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 746.75
$ws.Range("I18").Value = 329
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 329
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -45
$ws.Range("N18").Value = -2568
$ws.Range("H40").Value = 2399.1667
$ws.Range("I40").Value = 2398.75
$ws.Range("J40").Value = 2400
$ws.Range("K40").Value = 2398.75
$ws.Range("L40").Value = 2400
$ws.Range("M40").Value = -2223.75
$ws.Range("N40").Value = -2750
$ws.Range("H64").Value = 4306.087
$ws.Range("I64").Value = 3318.8
$ws.Range("J64").Value = 6157.25
$ws.Range("K64").Value = 3318.8
$ws.Range("L64").Value = 6157.25
$ws.Range("M64").Value = -3070.8
$ws.Range("N64").Value = -6653.25
$ws.Range("H67").Value = 4306.087
$ws.Range("I67").Value = 3318.8
$ws.Range("J67").Value = 6157.25
$ws.Range("K67").Value = 3318.8
$ws.Range("L67").Value = 6157.25
$ws.Range("M67").Value = -2460.8
$ws.Range("N67").Value = -7873.25
$ws.Range("H76").Value = 2270980.5
$ws.Range("I76").Value = 2927277.8
$ws.Range("J76").Value = 3772.5454
$ws.Range("K76").Value = 2927277.8
$ws.Range("L76").Value = 3772.5454
$ws.Range("M76").Value = -2926962.8
$ws.Range("N76").Value = -4402.5454
$ws.Range("H79").Value = 2270980.5
$ws.Range("I79").Value = 2927277.8
$ws.Range("J79").Value = 3772.5454
$ws.Range("K79").Value = 2927277.8
$ws.Range("L79").Value = 3772.5454
$ws.Range("M79").Value = -2926185.8
$ws.Range("N79").Value = -5956.5454
$ws.Range("H135").Value = 1716.1072
$ws.Range("I135").Value = 1721.8
$ws.Range("K135").Value = 15496.2
$ws.Range("M135").Value = -12961.2
$ws.Range("H138").Value = 8108382
$ws.Range("I138").Value = 2979552
$ws.Range("J138").Value = 10755520
$ws.Range("K138").Value = 8938656
$ws.Range("L138").Value = 32266560
$ws.Range("M138").Value = -8933516
$ws.Range("N138").Value = -32276840

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4243.7144
$ws.Range("I88").Value = 1581.2
$ws.Range("K88").Value = 1581.2
$ws.Range("M88").Value = -1175.2
$ws.Range("H91").Value = 4243.7144
$ws.Range("I91").Value = 1581.2
$ws.Range("K91").Value = 1581.2
$ws.Range("M91").Value = -177.2
$ws.Range("H110").Value = 1248.1111
$ws.Range("I110").Value = 1029.125
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1029.125
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 1015.875
$ws.Range("N110").Value = -7090

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8379.8125
$ws.Range("I86").Value = 2507
$ws.Range("K86").Value = 2507
$ws.Range("M86").Value = -1384
$ws.Range("H89").Value = 8379.8125
$ws.Range("I89").Value = 2507
$ws.Range("K89").Value = 12535
$ws.Range("M89").Value = -6919
$ws.Range("H94").Value = 1149.7916
$ws.Range("I94").Value = 1010.75
$ws.Range("J94").Value = 1845
$ws.Range("K94").Value = 1010.75
$ws.Range("L94").Value = 1845
$ws.Range("M94").Value = -559.75
$ws.Range("N94").Value = -2747
$ws.Range("H102").Value = 11052
$ws.Range("I102").Value = 11052
$ws.Range("K102").Value = 11052
$ws.Range("M102").Value = -7807
$ws.Range("H134").Value = 4441.85
$ws.Range("I134").Value = 3242.5833
$ws.Range("J134").Value = 6240.75
$ws.Range("K134").Value = 9727.749899999999
$ws.Range("L134").Value = 18722.25
$ws.Range("M134").Value = -7192.749899999999
$ws.Range("N134").Value = -23792.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1255.5
$ws.Range("I16").Value = 1340.6666
$ws.Range("K16").Value = 1340.6666
$ws.Range("M16").Value = -1053.6666
$ws.Range("H45").Value = 7066.7144
$ws.Range("I45").Value = 5067
$ws.Range("J45").Value = 7400
$ws.Range("K45").Value = 5067
$ws.Range("L45").Value = 7400
$ws.Range("M45").Value = -4474
$ws.Range("N45").Value = -8586
$ws.Range("H62").Value = 14683.444
$ws.Range("I62").Value = 21880.908
$ws.Range("J62").Value = 3373.1428
$ws.Range("K62").Value = 21880.908
$ws.Range("L62").Value = 3373.1428
$ws.Range("M62").Value = -21256.908
$ws.Range("N62").Value = -4621.1428
$ws.Range("H65").Value = 14683.444
$ws.Range("I65").Value = 21880.908
$ws.Range("J65").Value = 3373.1428
$ws.Range("K65").Value = 109404.54
$ws.Range("L65").Value = 16865.714
$ws.Range("M65").Value = -106284.54
$ws.Range("N65").Value = -23105.714
$ws.Range("H113").Value = 1255.5
$ws.Range("I113").Value = 1340.6666
$ws.Range("K113").Value = 1340.6666
$ws.Range("M113").Value = 829.3334
$ws.Range("H132").Value = 2849.6155
$ws.Range("I132").Value = 2538.0881
$ws.Range("K132").Value = 7614.2643
$ws.Range("M132").Value = -5084.2643

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1212.4375
$ws.Range("I5").Value = 547
$ws.Range("J5").Value = 2068
$ws.Range("K5").Value = 1641
$ws.Range("L5").Value = 6204
$ws.Range("M5").Value = -1529
$ws.Range("N5").Value = -6428
$ws.Range("H6").Value = 354.57144
$ws.Range("I6").Value = 96.40000000000001
$ws.Range("K6").Value = 289.2
$ws.Range("M6").Value = -176.2
$ws.Range("H26").Value = 204.75
$ws.Range("I26").Value = 191.5
$ws.Range("J26").Value = 271
$ws.Range("K26").Value = 574.5
$ws.Range("L26").Value = 813
$ws.Range("M26").Value = -286.5
$ws.Range("N26").Value = -1389
$ws.Range("H29").Value = 347.25
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -2662
$ws.Range("N41").Value = -3676
$ws.Range("H44").Value = 515.1667
$ws.Range("I44").Value = 198.5
$ws.Range("J44").Value = 673.5
$ws.Range("K44").Value = 595.5
$ws.Range("L44").Value = 2020.5
$ws.Range("M44").Value = -197.5
$ws.Range("N44").Value = -2816.5
$ws.Range("H55").Value = 1950
$ws.Range("I55").Value = 900
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 2700
$ws.Range("L55").Value = 9000
$ws.Range("M55").Value = -2523
$ws.Range("N55").Value = -9354
$ws.Range("H58").Value = 8193.538
$ws.Range("J58").Value = 9181.817999999999
$ws.Range("L58").Value = 27545.454
$ws.Range("N58").Value = -27801.454
$ws.Range("H64").Value = 1957.4615
$ws.Range("J64").Value = 5001
$ws.Range("L64").Value = 15003
$ws.Range("N64").Value = -15543
$ws.Range("H67").Value = 1957.4615
$ws.Range("J67").Value = 5001
$ws.Range("L67").Value = 15003
$ws.Range("N67").Value = -16875
$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 6000
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -8350
$ws.Range("H135").Value = 1212.4375
$ws.Range("I135").Value = 547
$ws.Range("J135").Value = 2068
$ws.Range("K135").Value = 4923
$ws.Range("L135").Value = 18612
$ws.Range("M135").Value = -2388
$ws.Range("N135").Value = -23682

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2308.5625
$ws.Range("I102").Value = 2217.6924
$ws.Range("K102").Value = 2217.6924
$ws.Range("M102").Value = -595.6923999999999
$ws.Range("H107").Value = 471.57144
$ws.Range("I107").Value = 432.66666
$ws.Range("J107").Value = 500.75
$ws.Range("K107").Value = 432.66666
$ws.Range("L107").Value = 500.75
$ws.Range("M107").Value = 1487.33334
$ws.Range("N107").Value = -4340.75
$ws.Range("H113").Value = 2245.5557
$ws.Range("I113").Value = 1003.3333
$ws.Range("J113").Value = 2866.6667
$ws.Range("K113").Value = 1003.3333
$ws.Range("L113").Value = 2866.6667
$ws.Range("M113").Value = 1166.6667
$ws.Range("N113").Value = -7206.6667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1055
$ws.Range("I46").Value = 729
$ws.Range("J46").Value = 1462.5
$ws.Range("K46").Value = 729
$ws.Range("L46").Value = 1462.5
$ws.Range("M46").Value = -541
$ws.Range("N46").Value = -1838.5
$ws.Range("H61").Value = 2837.5
$ws.Range("I61").Value = 2750
$ws.Range("J61").Value = 2925
$ws.Range("K61").Value = 2750
$ws.Range("L61").Value = 2925
$ws.Range("M61").Value = -2548
$ws.Range("N61").Value = -3329
$ws.Range("H68").Value = 1760.091
$ws.Range("I68").Value = 1774
$ws.Range("J68").Value = 1697.5
$ws.Range("K68").Value = 1774
$ws.Range("L68").Value = 1697.5
$ws.Range("M68").Value = -1025
$ws.Range("N68").Value = -3195.5
$ws.Range("H71").Value = 1760.091
$ws.Range("I71").Value = 1774
$ws.Range("J71").Value = 1697.5
$ws.Range("K71").Value = 8870
$ws.Range("L71").Value = 8487.5
$ws.Range("M71").Value = -5126
$ws.Range("N71").Value = -15975.5
$ws.Range("H109").Value = 24992.5
$ws.Range("J109").Value = 24992.5
$ws.Range("L109").Value = 24992.5
$ws.Range("N109").Value = -27766.5
$ws.Range("H113").Value = 2837.5
$ws.Range("I113").Value = 2750
$ws.Range("J113").Value = 2925
$ws.Range("K113").Value = 2750
$ws.Range("L113").Value = 2925
$ws.Range("M113").Value = -580
$ws.Range("N113").Value = -7265

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 450
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1350
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5190
$ws.Range("H122").Value = 144657.58
$ws.Range("I122").Value = 250900.75
$ws.Range("K122").Value = 752702.25
$ws.Range("M122").Value = -750252.25
